$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.264.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.330.85'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.10'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '185.72'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.325.27'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.181'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.579'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.40%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '47.11'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000268'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '666.35'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +10.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.859.92'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.49'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '66.339.82'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '17.94'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.327.40'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.13'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.897'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.79'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '102.08'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.99'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.54'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '31.74'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.69%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.50'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.81'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '597.02'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.30%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.88'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.00'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.106'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.848.06'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.64%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '56.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.70'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0₃0701'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.127'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '32.80'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.84%  '
$ws.Range('E43').Value = '  +5.61%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.17'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.29%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.338'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0414'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.03'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -11.29%  '
$ws.Range('E48').Value = '  -1.50%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.56'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.63%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  +1.82%  '
